$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 107.7653275
$ws.Range("H2").Value = 215.530655
$ws.Range("I2").Value = 0.1076850798785215
$ws.Range("J2").Value = 0.07714326402502852
$ws.Range("Q2").Value = 33.88472376937667
$ws.Range("R2").Value = 203.30834261626
$ws.Range("S2").Value = 0.1076850798785215
$ws.Range("T2").Value = 0.07714326402502852

# Row 3
$ws.Range("I3").Value = 0.2767710940483697
$ws.Range("J3").Value = 0.2974092456460348
$ws.Range("S3").Value = 0.2767710940483697
$ws.Range("T3").Value = 0.2974092456460348

# Row 4
$ws.Range("G4").Value = 210.8270723333333
$ws.Range("H4").Value = 632.481217
$ws.Range("I4").Value = 0.2106700796206445
$ws.Range("J4").Value = 0.2263792383218172
$ws.Range("Q4").Value = 66.29049690515157
$ws.Range("R4").Value = 596.6144721463641
$ws.Range("S4").Value = 0.2106700796206445
$ws.Range("T4").Value = 0.2263792383218172

# Row 5
$ws.Range("G5").Value = 100.56913
$ws.Range("H5").Value = 201.13826
$ws.Range("I5").Value = 0.1004942410383656
$ws.Range("J5").Value = 0.07199190248233985
$ws.Range("Q5").Value = 31.62201859198667
$ws.Range("R5").Value = 189.73211155192
$ws.Range("S5").Value = 0.1004942410383656
$ws.Range("T5").Value = 0.07199190248233985

# Row 6
$ws.Range("G6").Value = 102.5625483333333
$ws.Range("H6").Value = 307.687645
$ws.Range("I6").Value = 0.1024861749695859
$ws.Range("J6").Value = 0.1101283213539821
$ws.Range("Q6").Value = 32.24881044748223
$ws.Range("R6").Value = 290.23929402734
$ws.Range("S6").Value = 0.1024861749695859
$ws.Range("T6").Value = 0.1101283213539821

# Row 7
$ws.Range("G7").Value = 202.0437826666667
$ws.Range("H7").Value = 606.131348
$ws.Range("I7").Value = 0.2018933304445128
$ws.Range("J7").Value = 0.2169480281707975
$ws.Range("Q7").Value = 63.52876127973511
$ws.Range("R7").Value = 571.758851517616
$ws.Range("S7").Value = 0.2018933304445128
$ws.Range("T7").Value = 0.2169480281707975
